# Insert a new "corr0" worksheet before "corr1", matching the commit's
# addition of a correlated PARAM5/PARAM6 pair, and update the designinput
# sheet's row 9/10 to reference it plus bump B9 from 10 to 500.

$wb = $excel.ActiveWorkbook

# 1. Insert the new "corr0" sheet right before "corr1".
$ws_corr1 = $wb.Worksheets.Item("corr1")
$ws_corr0 = $wb.Worksheets.Add($ws_corr1)
$ws_corr0.Name = "corr0"

# Populate the correlation matrix for PARAM5 / PARAM6 (corr = 0.8)
$ws_corr0.Range("B1").Value = "PARAM5"
$ws_corr0.Range("C1").Value = "PARAM6"
$ws_corr0.Range("A2").Value = "PARAM5"
$ws_corr0.Range("B2").Value = 1
$ws_corr0.Range("A3").Value = "PARAM6"
$ws_corr0.Range("B3").Value = 0.8
$ws_corr0.Range("C3").Value = 1

# 2. Update designinput sheet: row 9 (PARAM5) and row 10 (PARAM6) now
#    point at the "corr0" correlation sheet, and numreal for PARAM5 bumps
#    from 10 to 500.
$ws_design = $wb.Worksheets.Item("designinput")
$ws_design.Range("B9").Value = 500
$ws_design.Range("O9").Value = "corr0"
$ws_design.Range("O10").Value = "corr0"

$ws_design.Range("B10").Select()

# 3. Finish with the new "corr0" sheet active/selected (matches the
#    saved workbook's activeTab + tabSelected state).
$ws_corr0.Activate()
$ws_corr0.Range("C8").Select()
